# Sheet1 held a flat shopping list (Client0 name + 4 loose items in column A).
# Reshape it into a 2-column table: Name | ShoppingList, with the items
# rolled up into a single comma-separated cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old A1:A5 list so nothing is left below the new 2x2 block.
$ws.Range("A1:A5").ClearContents()

# Header row.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "ShoppingList"

# Data row.
$ws.Range("A2").Value = "Client0"
$ws.Range("B2").Value = "Bread, Eggs, Milk, Apples, Oranges, Cucumber"

$ws.Range("B2").Select()
